$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header update
$ws.Range("A1").Value = "Datos actualizados a 18 de Octubre de 2020 a las 02:28"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 8341836
$ws.Range("C4").Value = 53403
$ws.Range("D4").Value = 5428995
$ws.Range("E4").Value = 2688563
$ws.Range("G4").Value = 634
$ws.Range("H4").Value = 224278

# Brasil (row 6)
$ws.Range("E6").Value = 435357
$ws.Range("G6").Value = 461
$ws.Range("H6").Value = 153690

# Paraguay (row 67)
$ws.Range("B67").Value = 54015
$ws.Range("C67").Value = 533
$ws.Range("D67").Value = 35524
$ws.Range("E67").Value = 17312
$ws.Range("G67").Value = 14
$ws.Range("H67").Value = 1179

# Noruega (row 96)
$ws.Range("B96").Value = 16369
$ws.Range("C96").Value = 97
$ws.Range("E96").Value = 4228

# Surinam (row 133)
$ws.Range("B133").Value = 5123
$ws.Range("C133").Value = 10
$ws.Range("D133").Value = 4936
$ws.Range("E133").Value = 78

# Uruguay overtakes Benin in the ranking (rows 155-156 swap country + row 155 gets
# Uruguay's fresh totals, row 156 keeps Benin's previous totals)
$ws.Range("A155").Value = "Uruguay"
$ws.Range("B155").Value = 2501
$ws.Range("C155").Value = 51
$ws.Range("D155").Value = 2052
$ws.Range("E155").Value = 398
$ws.Range("H155").Value = 51

$ws.Range("A156").Value = "Benin"
$ws.Range("B156").Value = 2496
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 2330
$ws.Range("E156").Value = 125
$ws.Range("H156").Value = 41

# San Martin (Parte Holandesa) (row 171)
$ws.Range("B171").Value = 749
$ws.Range("C171").Value = 3
$ws.Range("D171").Value = 664
$ws.Range("E171").Value = 63

# Islas Turcas y Caicos (row 174)
$ws.Range("B174").Value = 698
$ws.Range("C174").Value = 1
$ws.Range("D174").Value = 673
$ws.Range("E174").Value = 19
